# Automatic data update [2026-02-20 17:49]
# Refreshes extraction timestamps (col E) and several measured values
# (humidity, pressure, radiation, wind gust, temperatures) on the
# "Dades_Meteo" sheet to match the latest meteo.cat scrape.
#
# Note: a handful of humidity cells (col H) hold plain percent TEXT
# (e.g. "53%"), not numeric percentages. Excel's COM layer auto-parses a
# bare "NN%" string into a numeric percent value, so those assignments are
# prefixed with a leading apostrophe to force literal text entry, exactly
# like typing '53% directly into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-20 17:48:21'
$ws.Range('E3').Value = '2026-02-20 17:48:23'
$ws.Range('E4').Value = '2026-02-20 17:48:25'
$ws.Range('H4').Value = '''53%'
$ws.Range('J4').Value = '1022.0 hPa'
$ws.Range('E5').Value = '2026-02-20 17:48:28'
$ws.Range('E6').Value = '2026-02-20 17:48:30'
$ws.Range('J6').Value = '1022.0 hPa'
$ws.Range('O6').Value = '9.6 °C'
$ws.Range('E7').Value = '2026-02-20 17:48:33'
$ws.Range('H7').Value = '''44%'
$ws.Range('J7').Value = '1021.9 hPa'
$ws.Range('K7').Value = '12.1 MJ/m2'
$ws.Range('E8').Value = '2026-02-20 17:48:35'
$ws.Range('J8').Value = '1022.2 hPa'
$ws.Range('K8').Value = '11.0 MJ/m2'
$ws.Range('M8').Value = '13.2 °C 17:03 TU'
$ws.Range('O8').Value = '9.2 °C'
$ws.Range('E9').Value = '2026-02-20 17:48:38'
$ws.Range('E10').Value = '2026-02-20 17:48:40'
$ws.Range('E11').Value = '2026-02-20 17:48:41'
$ws.Range('E12').Value = '2026-02-20 17:48:42'
$ws.Range('E13').Value = '2026-02-20 17:48:44'
$ws.Range('J13').Value = '1022.8 hPa'
$ws.Range('L13').Value = '64.4 km/h - 5º 17:01 TU'
$ws.Range('E14').Value = '2026-02-20 17:48:45'
$ws.Range('K14').Value = '11.5 MJ/m2'
$ws.Range('O14').Value = '12.5 °C'
$ws.Range('E15').Value = '2026-02-20 17:48:46'
$ws.Range('E16').Value = '2026-02-20 17:48:47'
$ws.Range('O16').Value = '-3.7 °C'
$ws.Range('E17').Value = '2026-02-20 17:48:48'
$ws.Range('O17').Value = '2.7 °C'
$ws.Range('E18').Value = '2026-02-20 17:48:49'
$ws.Range('J18').Value = '1022.3 hPa'
$ws.Range('O18').Value = '8.1 °C'
$ws.Range('E19').Value = '2026-02-20 17:48:50'
$ws.Range('O19').Value = '4.3 °C'
$ws.Range('E20').Value = '2026-02-20 17:48:51'
$ws.Range('O20').Value = '-3.1 °C'
$ws.Range('E21').Value = '2026-02-20 17:48:52'
$ws.Range('H21').Value = '''35%'
$ws.Range('O21').Value = '9.5 °C'
$ws.Range('E22').Value = '2026-02-20 17:48:55'
$ws.Range('E23').Value = '2026-02-20 17:48:57'
$ws.Range('K23').Value = '16.0 MJ/m2'
$ws.Range('E24').Value = '2026-02-20 17:49:00'
$ws.Range('J24').Value = '1025.0 hPa'
$ws.Range('K24').Value = '14.0 MJ/m2'
$ws.Range('O24').Value = '9.5 °C'
$ws.Range('E25').Value = '2026-02-20 17:49:02'
$ws.Range('O25').Value = '-1.7 °C'
$ws.Range('E26').Value = '2026-02-20 17:49:05'
$ws.Range('O26').Value = '5.6 °C'
$ws.Range('E27').Value = '2026-02-20 17:49:07'
$ws.Range('O27').Value = '-0.8 °C'
$ws.Range('E28').Value = '2026-02-20 17:49:09'
$ws.Range('J28').Value = '1022.3 hPa'
$ws.Range('O28').Value = '7.3 °C'
$ws.Range('E29').Value = '2026-02-20 17:49:12'
$ws.Range('H29').Value = '''73%'
$ws.Range('E30').Value = '2026-02-20 17:49:14'
$ws.Range('J30').Value = '1021.6 hPa'
$ws.Range('E31').Value = '2026-02-20 17:49:17'
$ws.Range('J31').Value = '1020.8 hPa'
$ws.Range('E32').Value = '2026-02-20 17:49:19'
$ws.Range('O32').Value = '4.5 °C'
$ws.Range('E33').Value = '2026-02-20 17:49:21'
$ws.Range('J33').Value = '1022.2 hPa'
$ws.Range('E34').Value = '2026-02-20 17:49:24'
$ws.Range('L34').Value = '88.6 km/h - 16º 17:24 TU'
$ws.Range('M34').Value = '4.6 °C 17:25 TU'
$ws.Range('O34').Value = '0.3 °C'
$ws.Range('E35').Value = '2026-02-20 17:49:26'
$ws.Range('K35').Value = '10.9 MJ/m2'
$ws.Range('O35').Value = '3.9 °C'
$ws.Range('E36').Value = '2026-02-20 17:49:29'
$ws.Range('J36').Value = '1021.9 hPa'
$ws.Range('E37').Value = '2026-02-20 17:49:31'
$ws.Range('J37').Value = '1023.7 hPa'
$ws.Range('O37').Value = '5.1 °C'
$ws.Range('E38').Value = '2026-02-20 17:49:34'
$ws.Range('O38').Value = '9.0 °C'
$ws.Range('E39').Value = '2026-02-20 17:49:36'
$ws.Range('E40').Value = '2026-02-20 17:49:39'
$ws.Range('J40').Value = '1022.8 hPa'
$ws.Range('E41').Value = '2026-02-20 17:49:41'
$ws.Range('H41').Value = '''47%'
$ws.Range('J41').Value = '1022.6 hPa'
$ws.Range('O41').Value = '13.3 °C'
$ws.Range('E42').Value = '2026-02-20 17:49:44'
$ws.Range('H42').Value = '''68%'
$ws.Range('O42').Value = '9.9 °C'
$ws.Range('E43').Value = '2026-02-20 17:49:46'
$ws.Range('O43').Value = '4.9 °C'
$ws.Range('E44').Value = '2026-02-20 17:49:48'
$ws.Range('K44').Value = '10.0 MJ/m2'
$ws.Range('O44').Value = '-5.1 °C'
$ws.Range('E45').Value = '2026-02-20 17:49:51'
$ws.Range('J45').Value = '1029.1 hPa'
$ws.Range('K45').Value = '8.8 MJ/m2'
$ws.Range('E46').Value = '2026-02-20 17:49:53'
$ws.Range('J46').Value = '1026.1 hPa'
$ws.Range('K46').Value = '12.5 MJ/m2'
$ws.Range('O46').Value = '12.1 °C'
